# WORKING: multiple tabs Version 2
#
# Replaces the "drive letter" rows (C:/D: and their linked sites) with a
# single GOOGLE entry on the fieldnames/category/URL/color tabs, removing
# the now-unused rows 3-6 detail cells and their hyperlinks, and updates
# each sheet's remembered selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "fieldnames"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("fieldnames")

$ws.Range("A2").Value = "GOOGLE"

$ws.Range("A3:C3").ClearContents()
$ws.Range("C4:C6").ClearContents()

function Remove-HyperlinkAt($sheet, $addr) {
    foreach ($h in $sheet.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.Delete()
            return
        }
    }
}

Remove-HyperlinkAt $ws '$C$3'
Remove-HyperlinkAt $ws '$C$4'
Remove-HyperlinkAt $ws '$C$5'
Remove-HyperlinkAt $ws '$C$6'

$ws.Activate()
$ws.Range("A3").Select()

# ---------------------------------------------------------------------
# Sheet "category"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("category")

$ws.Range("A3:C3").ClearContents()
$ws.Range("C4:C6").ClearContents()

$ws.Activate()
$ws.Range("F11").Select()

# ---------------------------------------------------------------------
# Sheet "URL"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("URL")

$ws.Hyperlinks.Add($ws.Range("A2"), "http://www.google.de/", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "www.google.de")
$ws.Range("D2").Copy()
$ws.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A3:C3").ClearContents()
$ws.Range("C4:C6").ClearContents()

Remove-HyperlinkAt $ws '$B$3'

$ws.Activate()
$ws.Range("A3").Select()

# ---------------------------------------------------------------------
# Sheet "color"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("color")

$ws.Range("A3:C3").ClearContents()
$ws.Range("C4:C6").ClearContents()

$ws.Activate()
$ws.Range("G47").Select()
